$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. This shifts the old B,C,D (Start,
# Random, End) one column to the right, becoming C, D, E - matching the
# target layout where a new "Unnamed: 0" column is inserted right after
# the existing index column A.
$ws.Columns("B:B").Insert()

# New header cell for the inserted column.
$ws.Range("B1").Value = "Unnamed: 0"

# Give B1 the same header style (bold font, border, centered) as the
# other header cells by copying format from D1 (old "Random" header,
# shifted from C1, already carrying style index 1).
$ws.Range("D1").Copy()
$ws.Range("B1").PasteSpecial(-4122)

# The inserted column picked up column A's style (border/bold) for its
# data rows; the new "Unnamed: 0" data column should be plain/unstyled
# like the other data columns, so clear that formatting.
$ws.Range("B2:B11").ClearFormats()

# New column B holds the same integer sequence as the index column A.
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 3
$ws.Range("B6").Value = 4
$ws.Range("B7").Value = 5
$ws.Range("B8").Value = 6
$ws.Range("B9").Value = 7
$ws.Range("B10").Value = 8
$ws.Range("B11").Value = 9

# Updated timing values for the "Start" column (now column C).
$ws.Range("C2").Value = 0.000001452266666698657
$ws.Range("C3").Value = 0.000001647033333332123
$ws.Range("C4").Value = 0.000001945666666668128
$ws.Range("C5").Value = 0.000003070033333339477
$ws.Range("C6").Value = 0.000002372439999999188
$ws.Range("C7").Value = 0.000003595077777794762
$ws.Range("C8").Value = 0.00000289460000000136
$ws.Range("C9").Value = 0.00000357759166665043
$ws.Range("C10").Value = 0.000004196711111102626
$ws.Range("C11").Value = 0.000004456713333335453

# Updated timing values for the "Random" column (now column D).
$ws.Range("D2").Value = 0.000002357200000005833
$ws.Range("D3").Value = 0.000002271233333279573
$ws.Range("D4").Value = 0.000002401866666736573
$ws.Range("D5").Value = 0.000003259733333303908
$ws.Range("D6").Value = 0.0000026159200000014
$ws.Range("D7").Value = 0.000003491944444451495
$ws.Range("D8").Value = 0.000002910971428588839
$ws.Range("D9").Value = 0.000003141274999999647
$ws.Range("D10").Value = 0.00000372613333335418
$ws.Range("D11").Value = 0.000004456713333335453

# Updated timing values for the "End" column (now column E).
$ws.Range("E2").Value = 0.0000011311333331226100
$ws.Range("E3").Value = 0.000001145899999907366
$ws.Range("E4").Value = 0.000001155222222223529
$ws.Range("E5").Value = 0.000001501883333351846
$ws.Range("E6").Value = 0.000001152746666715151
$ws.Range("E7").Value = 0.000001375333333372611
$ws.Range("E8").Value = 0.000001159038095205948
$ws.Range("E9").Value = 0.000001198141666653403
$ws.Range("E10").Value = 0.000001292066666691324
$ws.Range("E11").Value = 0.000001439220000005055
